$d = $word.ActiveDocument

# 1) Date day stays, month 04 -> 06 (DE SOLICITUD month)
$d.Content.Find.Execute("04", $true, $false, $false, $false, $false, $true, 1, $false, "06", 2)

# 2) Year 1994 -> 2024
$d.Content.Find.Execute("1994", $true, $false, $false, $false, $false, $true, 1, $false, "2024", 2)

# 3) NORMAL / URGENTE checkbox toggle
$d.Content.Find.Execute("NORMAL:   ☐  URGENTE:  ⬛", $true, $false, $false, $false, $false, $true, 1, $false, "NORMAL:   ⬛  URGENTE:  ☐", 2)

# 4) Proyecto name -> multi-line with line breaks
$d.Content.Find.Execute("Proyecto Estadia", $true, $false, $false, $false, $false, $true, 1, $false, "Mantenimiento a^lla infraestructura^lfísica educativa^lde Educación^lSuperior", 2)

# 5) Actividad name -> multi-line with line breaks
$d.Content.Find.Execute("Pruebas de integración", $true, $false, $false, $false, $false, $true, 1, $false, "Verificación de la^lverificación del^lprograma anual^lde mantenimiento", 2)

# 6) Numeric quantity 34 -> 41
$d.Content.Find.Execute("34", $true, $false, $false, $false, $false, $true, 1, $false, "41", 2)

# 7) Paquete description text
$d.Content.Find.Execute("Asperiores numquam e", $true, $false, $false, $false, $false, $true, 1, $false, "Dolore est atque saepe neque temporibus dolores odio eaque nulla", 2)

# 8) Justificación text
$d.Content.Find.Execute("Magna sunt deleniti", $true, $false, $false, $false, $false, $true, 1, $false, "Asperiores tempore pariatur Sint voluptatem necessitatibus totam atque laborum Beatae ducimus similique", 2)
